# Add "PackageUrl" as a new column to the OSS-Inventory (ComponentsAndLicenses) sheet.
# The new column is inserted right after the "Usage pattern" column (column G),
# pushing all subsequent columns (old H..S) one position to the right (new I..T).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ComponentsAndLicenses")

# Insert a new, empty column at position H (8). Excel shifts the existing
# H:S columns to I:T and the new column inherits formatting from its left
# neighbour (column G).
$ws.Columns.Item(8).Insert()

# Match the new column's width to its neighbour, column G ("Usage pattern").
$ws.Columns.Item(8).ColumnWidth = $ws.Columns.Item(7).ColumnWidth()

# Populate the template placeholder (row 2) before the header text (row 1)
# so that the workbook's shared-string table records them in that order.
$ws.Cells.Item(2, 8).Value = "`$packageUrl`$"
$ws.Cells.Item(1, 8).Value = "PackageUrl"

# Update the sheet's active selection like a user would after making the edit.
$ws.Activate()
$ws.Range("G12").Select()
